$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2850
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2850
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2850
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3502
$ws.Range("H53").Value = 276.55554
$ws.Range("I53").Value = 110.333336
$ws.Range("J53").Value = 359.66666
$ws.Range("K53").Value = 110.333336
$ws.Range("L53").Value = 359.66666
$ws.Range("M53").Value = 526.666664
$ws.Range("N53").Value = -1633.66666
$ws.Range("H64").Value = 5608.2
$ws.Range("I64").Value = 4474.6
$ws.Range("J64").Value = 6741.8
$ws.Range("K64").Value = 4474.6
$ws.Range("L64").Value = 6741.8
$ws.Range("M64").Value = -4226.6
$ws.Range("N64").Value = -7237.8
$ws.Range("H67").Value = 5608.2
$ws.Range("I67").Value = 4474.6
$ws.Range("J67").Value = 6741.8
$ws.Range("K67").Value = 4474.6
$ws.Range("L67").Value = 6741.8
$ws.Range("M67").Value = -3616.6
$ws.Range("N67").Value = -8457.799999999999
$ws.Range("H111").Value = 15351.889
$ws.Range("I111").Value = 20694.5
$ws.Range("J111").Value = 4666.6665
$ws.Range("K111").Value = 62083.5
$ws.Range("L111").Value = 13999.9995
$ws.Range("M111").Value = -59016.5
$ws.Range("N111").Value = -20133.9995
$ws.Range("H132").Value = 973.7879
$ws.Range("I132").Value = 847.96875
$ws.Range("K132").Value = 2543.90625
$ws.Range("M132").Value = -13.90625
$ws.Range("H138").Value = 3169.3215
$ws.Range("I138").Value = 2615.6191
$ws.Range("J138").Value = 3353.889
$ws.Range("K138").Value = 7846.8573
$ws.Range("L138").Value = 10061.667
$ws.Range("M138").Value = -2706.8573
$ws.Range("N138").Value = -20341.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2092037.5
$ws.Range("J2").Value = 1139.8
$ws.Range("L2").Value = 1139.8
$ws.Range("N2").Value = -1365.8
$ws.Range("H32").Value = 35243.81
$ws.Range("I32").Value = 39396.35
$ws.Range("J32").Value = 24383.309
$ws.Range("K32").Value = 39396.35
$ws.Range("L32").Value = 24383.309
$ws.Range("M32").Value = -39109.35
$ws.Range("N32").Value = -24957.309
$ws.Range("H42").Value = 15999.5
$ws.Range("I42").Value = 15999
$ws.Range("J42").Value = 16000
$ws.Range("K42").Value = 15999
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = -15513
$ws.Range("N42").Value = -16972
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H74").Value = 235706.39
$ws.Range("I74").Value = 304566.75
$ws.Range("K74").Value = 304566.75
$ws.Range("M74").Value = -303692.75
$ws.Range("H77").Value = 235706.39
$ws.Range("I77").Value = 304566.75
$ws.Range("K77").Value = 1522833.75
$ws.Range("M77").Value = -1518465.75
$ws.Range("H116").Value = 2092037.5
$ws.Range("J116").Value = 1139.8
$ws.Range("L116").Value = 1139.8
$ws.Range("N116").Value = -5727.8
$ws.Range("H134").Value = 92630.75
$ws.Range("I134").Value = 50000
$ws.Range("J134").Value = 106841
$ws.Range("K134").Value = 50000
$ws.Range("L134").Value = 106841
$ws.Range("M134").Value = -44930
$ws.Range("N134").Value = -116981

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2092037.5
$ws.Range("J3").Value = 1139.8
$ws.Range("L3").Value = 1139.8
$ws.Range("N3").Value = -1367.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.73914
$ws.Range("I7").Value = 146.9375
$ws.Range("K7").Value = 146.9375
$ws.Range("M7").Value = -33.9375
$ws.Range("H22").Value = 513.9231
$ws.Range("I22").Value = 506.75
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 506.75
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -156.75
$ws.Range("N22").Value = -1300
$ws.Range("H31").Value = 32264248
$ws.Range("I31").Value = 90913760
$ws.Range("J31").Value = 7017.85
$ws.Range("K31").Value = 90913760
$ws.Range("L31").Value = 7017.85
$ws.Range("M31").Value = -90913465
$ws.Range("N31").Value = -7607.85
$ws.Range("H34").Value = 32264248
$ws.Range("I34").Value = 90913760
$ws.Range("J34").Value = 7017.85
$ws.Range("K34").Value = 90913760
$ws.Range("L34").Value = 7017.85
$ws.Range("M34").Value = -90913558
$ws.Range("N34").Value = -7421.85

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7338628.5
$ws.Range("I4").Value = 7505075
$ws.Range("K4").Value = 22515225
$ws.Range("M4").Value = -22515113
$ws.Range("H5").Value = 521.7857
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H17").Value = 32
$ws.Range("J17").Value = 32
$ws.Range("L17").Value = 96
$ws.Range("N17").Value = -434
$ws.Range("H55").Value = 2338
$ws.Range("J55").Value = 2264.4443
$ws.Range("L55").Value = 6793.3329
$ws.Range("N55").Value = -7147.3329
$ws.Range("H134").Value = 55557810
$ws.Range("I134").Value = 55557810
$ws.Range("K134").Value = 166673430
$ws.Range("M134").Value = -166668360
$ws.Range("H135").Value = 521.7857
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 75680
$ws.Range("J103").Value = 75680
$ws.Range("L103").Value = 75680
$ws.Range("N103").Value = -78024
$ws.Range("H123").Value = 43303.57
$ws.Range("J123").Value = 70833.336
$ws.Range("L123").Value = 70833.336
$ws.Range("N123").Value = -75733.336
$ws.Range("H134").Value = 107606.47
$ws.Range("J134").Value = 107606.47
$ws.Range("L134").Value = 322819.41
$ws.Range("N134").Value = -327889.41

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3185.3809
$ws.Range("J22").Value = 4021.6365
$ws.Range("L22").Value = 4021.6365
$ws.Range("N22").Value = -4611.636500000001
$ws.Range("H27").Value = 3185.3809
$ws.Range("J27").Value = 4021.6365
$ws.Range("L27").Value = 4021.6365
$ws.Range("N27").Value = -4235.636500000001
$ws.Range("H46").Value = 6910.4136
$ws.Range("J46").Value = 9005.143
$ws.Range("L46").Value = 9005.143
$ws.Range("N46").Value = -9381.143
$ws.Range("H68").Value = 5471.143
$ws.Range("I68").Value = 3833
$ws.Range("K68").Value = 3833
$ws.Range("M68").Value = -3084
$ws.Range("H71").Value = 5471.143
$ws.Range("I71").Value = 3833
$ws.Range("K71").Value = 19165
$ws.Range("M71").Value = -15421
$ws.Range("H132").Value = 5176.82
$ws.Range("I132").Value = 4510.483
$ws.Range("K132").Value = 13531.449
$ws.Range("M132").Value = -11001.449

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 14062.75
$ws.Range("I38").Value = 7083.6665
$ws.Range("K38").Value = 7083.6665
$ws.Range("M38").Value = -6610.6665
$ws.Range("H39").Value = 21666.666
$ws.Range("I39").Value = 15000
$ws.Range("J39").Value = 25000
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = -14587
$ws.Range("N39").Value = -25826
$ws.Range("H43").Value = 44665
$ws.Range("I43").Value = 56997.5
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 56997.5
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -56848.5
$ws.Range("N43").Value = -20298
$ws.Range("H132").Value = 6670.0586
$ws.Range("I132").Value = 2482.4443
$ws.Range("K132").Value = 7447.3329
$ws.Range("M132").Value = -4917.3329
